$d = $word.ActiveDocument

# Locate the paragraph that ends with the Dynamic Programming explanation
# ("...exponential to polynomial.") so the new content can be appended
# right after it, ahead of the existing trailing blank paragraphs.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*exponential to polynomial.*") {
        $target = $p
    }
}
if ($null -eq $target) {
    throw "Could not locate the Dynamic Programming paragraph to anchor the insertion."
}

# The new paragraphs to add, in order: a blank line, the "Remove Duplicate
# character" snippet (heading + python sample), then another blank line.
$lines = @(
    "",
    "Remove Duplicate character",
    'string="geeksforgeeks"',
    'p=""',
    "for char in string:",
    "    if char not in p:",
    "        p=p+char",
    "print(p)",
    ""
)

$cur = $target
foreach ($line in $lines) {
    $cur.Range.InsertParagraphAfter()
    $newPara = $cur.Next()
    if ($line -ne "") {
        $newPara.Range.InsertBefore($line)
    }
    $cur = $newPara
}

"done"
